$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new as-of date (Feb 11, 2022)
$ws.Name = "Through 2022-02-11"

# Update the running-month header label (B1); the rest of row 1 (historical
# month headers C1:Q1) is untouched by this commit.
$ws.Range("B1").Value = 'February 2022 (through February 11)'

# Refresh the neighborhood rows (A2:Q86): new Feb-19 data changes the running
# Feb-2022 totals, which re-sorts the whole table (it is sorted by the fixed
# January-2022 column, with ties broken by the groupby/sort order coming out
# of the underlying pipeline), so many rows pick up a different neighborhood
# and/or different monthly counts.
$data = New-Object 'object[,]' 85,17
$data[0,0] = 'Englewood'
$data[0,1] = 3
$data[0,2] = 12
$data[0,3] = $null
$data[0,4] = 4
$data[0,5] = $null
$data[0,6] = 5
$data[0,7] = 1
$data[0,8] = 5
$data[0,9] = $null
$data[0,10] = 2
$data[0,11] = 2
$data[0,12] = $null
$data[0,13] = 1
$data[0,14] = 5
$data[0,15] = $null
$data[0,16] = 3
$data[1,0] = 'Austin'
$data[1,1] = 2
$data[1,2] = 11
$data[1,3] = 6
$data[1,4] = 10
$data[1,5] = 1
$data[1,6] = 7
$data[1,7] = 1
$data[1,8] = 3
$data[1,9] = 3
$data[1,10] = 8
$data[1,11] = 4
$data[1,12] = 11
$data[1,13] = $null
$data[1,14] = 3
$data[1,15] = $null
$data[1,16] = 4
$data[2,0] = 'New City'
$data[2,1] = $null
$data[2,2] = 10
$data[2,3] = 1
$data[2,4] = 3
$data[2,5] = $null
$data[2,6] = 1
$data[2,7] = $null
$data[2,8] = 1
$data[2,9] = $null
$data[2,10] = $null
$data[2,11] = 1
$data[2,12] = 1
$data[2,13] = 1
$data[2,14] = $null
$data[2,15] = $null
$data[2,16] = $null
$data[3,0] = 'Calumet Heights'
$data[3,1] = 2
$data[3,2] = 6
$data[3,3] = $null
$data[3,4] = 17
$data[3,5] = $null
$data[3,6] = $null
$data[3,7] = $null
$data[3,8] = 1
$data[3,9] = $null
$data[3,10] = 1
$data[3,11] = $null
$data[3,12] = 1
$data[3,13] = $null
$data[3,14] = $null
$data[3,15] = $null
$data[3,16] = $null
$data[4,0] = 'South Shore'
$data[4,1] = 3
$data[4,2] = 6
$data[4,3] = 7
$data[4,4] = 9
$data[4,5] = 1
$data[4,6] = 2
$data[4,7] = $null
$data[4,8] = $null
$data[4,9] = $null
$data[4,10] = $null
$data[4,11] = 2
$data[4,12] = 4
$data[4,13] = $null
$data[4,14] = 1
$data[4,15] = $null
$data[4,16] = 1
$data[5,0] = 'Auburn Gresham'
$data[5,1] = $null
$data[5,2] = 6
$data[5,3] = 1
$data[5,4] = 4
$data[5,5] = $null
$data[5,6] = $null
$data[5,7] = $null
$data[5,8] = 3
$data[5,9] = 1
$data[5,10] = 1
$data[5,11] = $null
$data[5,12] = 3
$data[5,13] = 1
$data[5,14] = 2
$data[5,15] = $null
$data[5,16] = 1
$data[6,0] = 'North Lawndale'
$data[6,1] = 2
$data[6,2] = 6
$data[6,3] = 3
$data[6,4] = 10
$data[6,5] = 3
$data[6,6] = 4
$data[6,7] = $null
$data[6,8] = 1
$data[6,9] = 1
$data[6,10] = 2
$data[6,11] = $null
$data[6,12] = 4
$data[6,13] = $null
$data[6,14] = 4
$data[6,15] = $null
$data[6,16] = 3
$data[7,0] = 'United Center'
$data[7,1] = 1
$data[7,2] = 5
$data[7,3] = 2
$data[7,4] = 4
$data[7,5] = $null
$data[7,6] = $null
$data[7,7] = $null
$data[7,8] = 3
$data[7,9] = 1
$data[7,10] = 2
$data[7,11] = $null
$data[7,12] = 2
$data[7,13] = $null
$data[7,14] = 2
$data[7,15] = $null
$data[7,16] = $null
$data[8,0] = 'Grand Crossing'
$data[8,1] = 1
$data[8,2] = 5
$data[8,3] = 1
$data[8,4] = 4
$data[8,5] = $null
$data[8,6] = 1
$data[8,7] = 2
$data[8,8] = 1
$data[8,9] = 1
$data[8,10] = 1
$data[8,11] = 2
$data[8,12] = 3
$data[8,13] = $null
$data[8,14] = 2
$data[8,15] = $null
$data[8,16] = 1
$data[9,0] = 'Edgewater'
$data[9,1] = $null
$data[9,2] = 4
$data[9,3] = 1
$data[9,4] = 4
$data[9,5] = $null
$data[9,6] = $null
$data[9,7] = 1
$data[9,8] = 1
$data[9,9] = $null
$data[9,10] = $null
$data[9,11] = $null
$data[9,12] = $null
$data[9,13] = $null
$data[9,14] = $null
$data[9,15] = $null
$data[9,16] = $null
$data[10,0] = 'Little Italy, UIC'
$data[10,1] = 1
$data[10,2] = 4
$data[10,3] = 4
$data[10,4] = 1
$data[10,5] = $null
$data[10,6] = $null
$data[10,7] = $null
$data[10,8] = 2
$data[10,9] = 2
$data[10,10] = 1
$data[10,11] = $null
$data[10,12] = 1
$data[10,13] = $null
$data[10,14] = 1
$data[10,15] = 1
$data[10,16] = 2
$data[11,0] = 'Bridgeport'
$data[11,1] = 1
$data[11,2] = 4
$data[11,3] = 1
$data[11,4] = 2
$data[11,5] = $null
$data[11,6] = $null
$data[11,7] = $null
$data[11,8] = $null
$data[11,9] = $null
$data[11,10] = $null
$data[11,11] = $null
$data[11,12] = 1
$data[11,13] = $null
$data[11,14] = $null
$data[11,15] = $null
$data[11,16] = $null
$data[12,0] = 'West Town'
$data[12,1] = 2
$data[12,2] = 4
$data[12,3] = 4
$data[12,4] = 1
$data[12,5] = 1
$data[12,6] = 1
$data[12,7] = 1
$data[12,8] = $null
$data[12,9] = $null
$data[12,10] = $null
$data[12,11] = 1
$data[12,12] = 2
$data[12,13] = $null
$data[12,14] = $null
$data[12,15] = $null
$data[12,16] = $null
$data[13,0] = 'Garfield Park'
$data[13,1] = 4
$data[13,2] = 4
$data[13,3] = 5
$data[13,4] = 12
$data[13,5] = 4
$data[13,6] = 10
$data[13,7] = 1
$data[13,8] = 2
$data[13,9] = 1
$data[13,10] = 6
$data[13,11] = 1
$data[13,12] = 5
$data[13,13] = 1
$data[13,14] = 6
$data[13,15] = $null
$data[13,16] = 1
$data[14,0] = 'Kenwood'
$data[14,1] = 1
$data[14,2] = 3
$data[14,3] = 1
$data[14,4] = 8
$data[14,5] = $null
$data[14,6] = 1
$data[14,7] = $null
$data[14,8] = $null
$data[14,9] = $null
$data[14,10] = 2
$data[14,11] = $null
$data[14,12] = $null
$data[14,13] = $null
$data[14,14] = $null
$data[14,15] = $null
$data[14,16] = $null
$data[15,0] = 'Uptown'
$data[15,1] = $null
$data[15,2] = 3
$data[15,3] = 1
$data[15,4] = 1
$data[15,5] = $null
$data[15,6] = $null
$data[15,7] = $null
$data[15,8] = $null
$data[15,9] = $null
$data[15,10] = $null
$data[15,11] = $null
$data[15,12] = $null
$data[15,13] = $null
$data[15,14] = 1
$data[15,15] = $null
$data[15,16] = $null
$data[16,0] = 'Humboldt Park'
$data[16,1] = $null
$data[16,2] = 3
$data[16,3] = $null
$data[16,4] = 7
$data[16,5] = 1
$data[16,6] = 5
$data[16,7] = $null
$data[16,8] = 5
$data[16,9] = 3
$data[16,10] = 6
$data[16,11] = 2
$data[16,12] = 5
$data[16,13] = $null
$data[16,14] = 4
$data[16,15] = $null
$data[16,16] = $null
$data[17,0] = 'Logan Square'
$data[17,1] = $null
$data[17,2] = 3
$data[17,3] = $null
$data[17,4] = 1
$data[17,5] = $null
$data[17,6] = 1
$data[17,7] = $null
$data[17,8] = 2
$data[17,9] = 1
$data[17,10] = 2
$data[17,11] = $null
$data[17,12] = 3
$data[17,13] = $null
$data[17,14] = 2
$data[17,15] = $null
$data[17,16] = $null
$data[18,0] = 'Grand Boulevard'
$data[18,1] = $null
$data[18,2] = 3
$data[18,3] = 2
$data[18,4] = 9
$data[18,5] = $null
$data[18,6] = 1
$data[18,7] = $null
$data[18,8] = $null
$data[18,9] = $null
$data[18,10] = 2
$data[18,11] = $null
$data[18,12] = $null
$data[18,13] = $null
$data[18,14] = $null
$data[18,15] = 1
$data[18,16] = $null
$data[19,0] = 'Archer Heights'
$data[19,1] = 1
$data[19,2] = 3
$data[19,3] = $null
$data[19,4] = 3
$data[19,5] = $null
$data[19,6] = $null
$data[19,7] = $null
$data[19,8] = $null
$data[19,9] = $null
$data[19,10] = $null
$data[19,11] = $null
$data[19,12] = $null
$data[19,13] = $null
$data[19,14] = 1
$data[19,15] = $null
$data[19,16] = $null
$data[20,0] = 'Chicago Lawn'
$data[20,1] = 4
$data[20,2] = 3
$data[20,3] = $null
$data[20,4] = 6
$data[20,5] = $null
$data[20,6] = 3
$data[20,7] = $null
$data[20,8] = 1
$data[20,9] = $null
$data[20,10] = 4
$data[20,11] = $null
$data[20,12] = $null
$data[20,13] = 3
$data[20,14] = $null
$data[20,15] = $null
$data[20,16] = $null
$data[21,0] = 'Chatham'
$data[21,1] = 1
$data[21,2] = 3
$data[21,3] = $null
$data[21,4] = 3
$data[21,5] = 1
$data[21,6] = 4
$data[21,7] = 1
$data[21,8] = $null
$data[21,9] = $null
$data[21,10] = 1
$data[21,11] = $null
$data[21,12] = $null
$data[21,13] = 1
$data[21,14] = $null
$data[21,15] = $null
$data[21,16] = 1
$data[22,0] = 'Near South Side'
$data[22,1] = 1
$data[22,2] = 3
$data[22,3] = $null
$data[22,4] = 2
$data[22,5] = $null
$data[22,6] = $null
$data[22,7] = $null
$data[22,8] = 1
$data[22,9] = $null
$data[22,10] = $null
$data[22,11] = $null
$data[22,12] = 1
$data[22,13] = $null
$data[22,14] = $null
$data[22,15] = $null
$data[22,16] = $null
$data[23,0] = 'Morgan Park'
$data[23,1] = $null
$data[23,2] = 2
$data[23,3] = $null
$data[23,4] = 1
$data[23,5] = $null
$data[23,6] = 1
$data[23,7] = $null
$data[23,8] = $null
$data[23,9] = $null
$data[23,10] = 1
$data[23,11] = 2
$data[23,12] = 1
$data[23,13] = $null
$data[23,14] = $null
$data[23,15] = $null
$data[23,16] = $null
$data[24,0] = 'South Deering'
$data[24,1] = $null
$data[24,2] = 2
$data[24,3] = $null
$data[24,4] = 1
$data[24,5] = $null
$data[24,6] = $null
$data[24,7] = $null
$data[24,8] = 1
$data[24,9] = $null
$data[24,10] = $null
$data[24,11] = $null
$data[24,12] = $null
$data[24,13] = $null
$data[24,14] = $null
$data[24,15] = $null
$data[24,16] = $null
$data[25,0] = 'Roseland'
$data[25,1] = $null
$data[25,2] = 2
$data[25,3] = $null
$data[25,4] = 9
$data[25,5] = 1
$data[25,6] = 6
$data[25,7] = $null
$data[25,8] = 1
$data[25,9] = $null
$data[25,10] = 3
$data[25,11] = 1
$data[25,12] = 6
$data[25,13] = $null
$data[25,14] = 2
$data[25,15] = $null
$data[25,16] = 1
$data[26,0] = 'North Park'
$data[26,1] = 1
$data[26,2] = 2
$data[26,3] = $null
$data[26,4] = $null
$data[26,5] = 1
$data[26,6] = $null
$data[26,7] = $null
$data[26,8] = 1
$data[26,9] = $null
$data[26,10] = $null
$data[26,11] = $null
$data[26,12] = $null
$data[26,13] = $null
$data[26,14] = $null
$data[26,15] = $null
$data[26,16] = $null
$data[27,0] = 'Washington Heights'
$data[27,1] = 1
$data[27,2] = 2
$data[27,3] = $null
$data[27,4] = 4
$data[27,5] = $null
$data[27,6] = 1
$data[27,7] = $null
$data[27,8] = 2
$data[27,9] = $null
$data[27,10] = $null
$data[27,11] = $null
$data[27,12] = 2
$data[27,13] = 1
$data[27,14] = 2
$data[27,15] = $null
$data[27,16] = $null
$data[28,0] = 'Portage Park'
$data[28,1] = $null
$data[28,2] = 2
$data[28,3] = $null
$data[28,4] = 3
$data[28,5] = $null
$data[28,6] = 1
$data[28,7] = $null
$data[28,8] = $null
$data[28,9] = 1
$data[28,10] = 3
$data[28,11] = $null
$data[28,12] = $null
$data[28,13] = $null
$data[28,14] = $null
$data[28,15] = $null
$data[28,16] = $null
$data[29,0] = 'Hyde Park'
$data[29,1] = $null
$data[29,2] = 2
$data[29,3] = $null
$data[29,4] = 4
$data[29,5] = $null
$data[29,6] = 1
$data[29,7] = $null
$data[29,8] = $null
$data[29,9] = $null
$data[29,10] = 3
$data[29,11] = $null
$data[29,12] = $null
$data[29,13] = $null
$data[29,14] = $null
$data[29,15] = $null
$data[29,16] = $null
$data[30,0] = 'Riverdale'
$data[30,1] = $null
$data[30,2] = 2
$data[30,3] = $null
$data[30,4] = $null
$data[30,5] = $null
$data[30,6] = $null
$data[30,7] = $null
$data[30,8] = $null
$data[30,9] = $null
$data[30,10] = 1
$data[30,11] = $null
$data[30,12] = $null
$data[30,13] = $null
$data[30,14] = $null
$data[30,15] = $null
$data[30,16] = $null
$data[31,0] = 'West Lawn'
$data[31,1] = 2
$data[31,2] = 2
$data[31,3] = $null
$data[31,4] = 1
$data[31,5] = $null
$data[31,6] = $null
$data[31,7] = $null
$data[31,8] = $null
$data[31,9] = $null
$data[31,10] = $null
$data[31,11] = $null
$data[31,12] = 1
$data[31,13] = $null
$data[31,14] = 1
$data[31,15] = $null
$data[31,16] = $null
$data[32,0] = 'West Loop'
$data[32,1] = 3
$data[32,2] = 2
$data[32,3] = $null
$data[32,4] = 2
$data[32,5] = $null
$data[32,6] = 3
$data[32,7] = $null
$data[32,8] = 1
$data[32,9] = 1
$data[32,10] = 2
$data[32,11] = $null
$data[32,12] = 1
$data[32,13] = $null
$data[32,14] = 2
$data[32,15] = $null
$data[32,16] = $null
$data[33,0] = 'West Ridge'
$data[33,1] = 1
$data[33,2] = 2
$data[33,3] = 1
$data[33,4] = 2
$data[33,5] = 3
$data[33,6] = 1
$data[33,7] = $null
$data[33,8] = 2
$data[33,9] = $null
$data[33,10] = 2
$data[33,11] = $null
$data[33,12] = $null
$data[33,13] = $null
$data[33,14] = 2
$data[33,15] = $null
$data[33,16] = 1
$data[34,0] = 'South Chicago'
$data[34,1] = $null
$data[34,2] = 1
$data[34,3] = $null
$data[34,4] = 3
$data[34,5] = $null
$data[34,6] = $null
$data[34,7] = $null
$data[34,8] = $null
$data[34,9] = 1
$data[34,10] = 1
$data[34,11] = $null
$data[34,12] = $null
$data[34,13] = $null
$data[34,14] = $null
$data[34,15] = 1
$data[34,16] = $null
$data[35,0] = 'West Pullman'
$data[35,1] = 1
$data[35,2] = 1
$data[35,3] = $null
$data[35,4] = $null
$data[35,5] = $null
$data[35,6] = 2
$data[35,7] = $null
$data[35,8] = $null
$data[35,9] = $null
$data[35,10] = 1
$data[35,11] = $null
$data[35,12] = $null
$data[35,13] = $null
$data[35,14] = $null
$data[35,15] = $null
$data[35,16] = $null
$data[36,0] = 'Wicker Park'
$data[36,1] = $null
$data[36,2] = 1
$data[36,3] = 2
$data[36,4] = 7
$data[36,5] = 1
$data[36,6] = 2
$data[36,7] = $null
$data[36,8] = $null
$data[36,9] = $null
$data[36,10] = $null
$data[36,11] = 1
$data[36,12] = $null
$data[36,13] = $null
$data[36,14] = $null
$data[36,15] = $null
$data[36,16] = $null
$data[37,0] = 'Rogers Park'
$data[37,1] = $null
$data[37,2] = 1
$data[37,3] = $null
$data[37,4] = 2
$data[37,5] = $null
$data[37,6] = $null
$data[37,7] = $null
$data[37,8] = 1
$data[37,9] = $null
$data[37,10] = 1
$data[37,11] = $null
$data[37,12] = $null
$data[37,13] = $null
$data[37,14] = $null
$data[37,15] = $null
$data[37,16] = 1
$data[38,0] = 'Albany Park'
$data[38,1] = 1
$data[38,2] = 1
$data[38,3] = 1
$data[38,4] = 3
$data[38,5] = $null
$data[38,6] = $null
$data[38,7] = $null
$data[38,8] = $null
$data[38,9] = $null
$data[38,10] = 1
$data[38,11] = 1
$data[38,12] = 1
$data[38,13] = $null
$data[38,14] = 1
$data[38,15] = $null
$data[38,16] = $null
$data[39,0] = 'Loop'
$data[39,1] = $null
$data[39,2] = 1
$data[39,3] = 1
$data[39,4] = 1
$data[39,5] = $null
$data[39,6] = 1
$data[39,7] = $null
$data[39,8] = $null
$data[39,9] = $null
$data[39,10] = 1
$data[39,11] = $null
$data[39,12] = $null
$data[39,13] = $null
$data[39,14] = $null
$data[39,15] = $null
$data[39,16] = $null
$data[40,0] = 'River North'
$data[40,1] = 1
$data[40,2] = 1
$data[40,3] = $null
$data[40,4] = 3
$data[40,5] = $null
$data[40,6] = 1
$data[40,7] = $null
$data[40,8] = 1
$data[40,9] = $null
$data[40,10] = $null
$data[40,11] = $null
$data[40,12] = $null
$data[40,13] = $null
$data[40,14] = $null
$data[40,15] = $null
$data[40,16] = $null
$data[41,0] = 'Old Town'
$data[41,1] = $null
$data[41,2] = 1
$data[41,3] = $null
$data[41,4] = 2
$data[41,5] = 1
$data[41,6] = $null
$data[41,7] = $null
$data[41,8] = $null
$data[41,9] = $null
$data[41,10] = $null
$data[41,11] = $null
$data[41,12] = $null
$data[41,13] = $null
$data[41,14] = $null
$data[41,15] = $null
$data[41,16] = $null
$data[42,0] = 'Mount Greenwood'
$data[42,1] = $null
$data[42,2] = 1
$data[42,3] = $null
$data[42,4] = 1
$data[42,5] = $null
$data[42,6] = $null
$data[42,7] = $null
$data[42,8] = $null
$data[42,9] = $null
$data[42,10] = $null
$data[42,11] = $null
$data[42,12] = $null
$data[42,13] = $null
$data[42,14] = $null
$data[42,15] = $null
$data[42,16] = $null
$data[43,0] = 'Mckinley Park'
$data[43,1] = $null
$data[43,2] = 1
$data[43,3] = 1
$data[43,4] = $null
$data[43,5] = $null
$data[43,6] = $null
$data[43,7] = $null
$data[43,8] = 1
$data[43,9] = $null
$data[43,10] = $null
$data[43,11] = $null
$data[43,12] = $null
$data[43,13] = $null
$data[43,14] = $null
$data[43,15] = $null
$data[43,16] = $null
$data[44,0] = 'Lower West Side'
$data[44,1] = $null
$data[44,2] = 1
$data[44,3] = $null
$data[44,4] = 2
$data[44,5] = $null
$data[44,6] = $null
$data[44,7] = $null
$data[44,8] = 1
$data[44,9] = $null
$data[44,10] = $null
$data[44,11] = $null
$data[44,12] = $null
$data[44,13] = $null
$data[44,14] = $null
$data[44,15] = $null
$data[44,16] = $null
$data[45,0] = 'Little Village'
$data[45,1] = 1
$data[45,2] = 1
$data[45,3] = $null
$data[45,4] = 2
$data[45,5] = $null
$data[45,6] = $null
$data[45,7] = 1
$data[45,8] = $null
$data[45,9] = $null
$data[45,10] = 2
$data[45,11] = 1
$data[45,12] = $null
$data[45,13] = $null
$data[45,14] = $null
$data[45,15] = $null
$data[45,16] = $null
$data[46,0] = 'Lake View'
$data[46,1] = 1
$data[46,2] = 1
$data[46,3] = $null
$data[46,4] = 2
$data[46,5] = $null
$data[46,6] = $null
$data[46,7] = $null
$data[46,8] = $null
$data[46,9] = $null
$data[46,10] = $null
$data[46,11] = $null
$data[46,12] = $null
$data[46,13] = $null
$data[46,14] = 1
$data[46,15] = $null
$data[46,16] = 1
$data[47,0] = 'Irving Park'
$data[47,1] = $null
$data[47,2] = 1
$data[47,3] = $null
$data[47,4] = 1
$data[47,5] = $null
$data[47,6] = $null
$data[47,7] = $null
$data[47,8] = $null
$data[47,9] = $null
$data[47,10] = 3
$data[47,11] = $null
$data[47,12] = 3
$data[47,13] = $null
$data[47,14] = $null
$data[47,15] = $null
$data[47,16] = 1
$data[48,0] = 'Hermosa'
$data[48,1] = 1
$data[48,2] = 1
$data[48,3] = $null
$data[48,4] = $null
$data[48,5] = $null
$data[48,6] = $null
$data[48,7] = $null
$data[48,8] = $null
$data[48,9] = $null
$data[48,10] = $null
$data[48,11] = 1
$data[48,12] = $null
$data[48,13] = $null
$data[48,14] = $null
$data[48,15] = $null
$data[48,16] = $null
$data[49,0] = 'East Side'
$data[49,1] = $null
$data[49,2] = 1
$data[49,3] = $null
$data[49,4] = $null
$data[49,5] = $null
$data[49,6] = $null
$data[49,7] = $null
$data[49,8] = $null
$data[49,9] = $null
$data[49,10] = $null
$data[49,11] = $null
$data[49,12] = $null
$data[49,13] = $null
$data[49,14] = $null
$data[49,15] = $null
$data[49,16] = $null
$data[50,0] = 'Bucktown'
$data[50,1] = $null
$data[50,2] = 1
$data[50,3] = 1
$data[50,4] = 2
$data[50,5] = $null
$data[50,6] = $null
$data[50,7] = $null
$data[50,8] = $null
$data[50,9] = $null
$data[50,10] = 1
$data[50,11] = $null
$data[50,12] = 1
$data[50,13] = $null
$data[50,14] = $null
$data[50,15] = $null
$data[50,16] = $null
$data[51,0] = 'Brighton Park'
$data[51,1] = 1
$data[51,2] = 1
$data[51,3] = $null
$data[51,4] = 1
$data[51,5] = $null
$data[51,6] = $null
$data[51,7] = $null
$data[51,8] = $null
$data[51,9] = $null
$data[51,10] = $null
$data[51,11] = $null
$data[51,12] = $null
$data[51,13] = 1
$data[51,14] = $null
$data[51,15] = $null
$data[51,16] = $null
$data[52,0] = 'Belmont Cragin'
$data[52,1] = $null
$data[52,2] = 1
$data[52,3] = 2
$data[52,4] = $null
$data[52,5] = $null
$data[52,6] = 1
$data[52,7] = $null
$data[52,8] = 2
$data[52,9] = $null
$data[52,10] = 2
$data[52,11] = 2
$data[52,12] = 2
$data[52,13] = $null
$data[52,14] = $null
$data[52,15] = $null
$data[52,16] = 1
$data[53,0] = 'Avalon Park'
$data[53,1] = $null
$data[53,2] = 1
$data[53,3] = 1
$data[53,4] = 3
$data[53,5] = 1
$data[53,6] = $null
$data[53,7] = $null
$data[53,8] = $null
$data[53,9] = $null
$data[53,10] = 1
$data[53,11] = 1
$data[53,12] = $null
$data[53,13] = $null
$data[53,14] = $null
$data[53,15] = $null
$data[53,16] = $null
$data[54,0] = 'Ashburn'
$data[54,1] = $null
$data[54,2] = 1
$data[54,3] = $null
$data[54,4] = $null
$data[54,5] = $null
$data[54,6] = $null
$data[54,7] = $null
$data[54,8] = $null
$data[54,9] = $null
$data[54,10] = 3
$data[54,11] = $null
$data[54,12] = $null
$data[54,13] = $null
$data[54,14] = 1
$data[54,15] = $null
$data[54,16] = $null
$data[55,0] = 'Woodlawn'
$data[55,1] = 4
$data[55,2] = 1
$data[55,3] = $null
$data[55,4] = 3
$data[55,5] = 1
$data[55,6] = $null
$data[55,7] = $null
$data[55,8] = $null
$data[55,9] = $null
$data[55,10] = 1
$data[55,11] = $null
$data[55,12] = 1
$data[55,13] = $null
$data[55,14] = $null
$data[55,15] = $null
$data[55,16] = 2
$data[56,0] = 'Andersonville'
$data[56,1] = $null
$data[56,2] = $null
$data[56,3] = $null
$data[56,4] = 1
$data[56,5] = $null
$data[56,6] = $null
$data[56,7] = $null
$data[56,8] = $null
$data[56,9] = $null
$data[56,10] = $null
$data[56,11] = $null
$data[56,12] = $null
$data[56,13] = $null
$data[56,14] = $null
$data[56,15] = $null
$data[56,16] = $null
$data[57,0] = 'Avondale'
$data[57,1] = $null
$data[57,2] = $null
$data[57,3] = 1
$data[57,4] = 2
$data[57,5] = $null
$data[57,6] = $null
$data[57,7] = $null
$data[57,8] = $null
$data[57,9] = $null
$data[57,10] = 1
$data[57,11] = $null
$data[57,12] = 1
$data[57,13] = $null
$data[57,14] = $null
$data[57,15] = $null
$data[57,16] = $null
$data[58,0] = 'Chinatown'
$data[58,1] = 1
$data[58,2] = $null
$data[58,3] = $null
$data[58,4] = 3
$data[58,5] = $null
$data[58,6] = $null
$data[58,7] = $null
$data[58,8] = $null
$data[58,9] = $null
$data[58,10] = $null
$data[58,11] = $null
$data[58,12] = $null
$data[58,13] = $null
$data[58,14] = $null
$data[58,15] = $null
$data[58,16] = $null
$data[59,0] = 'Clearing'
$data[59,1] = $null
$data[59,2] = $null
$data[59,3] = 1
$data[59,4] = $null
$data[59,5] = $null
$data[59,6] = $null
$data[59,7] = $null
$data[59,8] = $null
$data[59,9] = $null
$data[59,10] = $null
$data[59,11] = $null
$data[59,12] = $null
$data[59,13] = $null
$data[59,14] = $null
$data[59,15] = $null
$data[59,16] = 1
$data[60,0] = 'Douglas'
$data[60,1] = $null
$data[60,2] = $null
$data[60,3] = $null
$data[60,4] = 3
$data[60,5] = $null
$data[60,6] = 1
$data[60,7] = $null
$data[60,8] = 1
$data[60,9] = $null
$data[60,10] = 3
$data[60,11] = $null
$data[60,12] = $null
$data[60,13] = $null
$data[60,14] = $null
$data[60,15] = $null
$data[60,16] = $null
$data[61,0] = 'Dunning'
$data[61,1] = $null
$data[61,2] = $null
$data[61,3] = $null
$data[61,4] = $null
$data[61,5] = $null
$data[61,6] = $null
$data[61,7] = $null
$data[61,8] = $null
$data[61,9] = 1
$data[61,10] = $null
$data[61,11] = $null
$data[61,12] = $null
$data[61,13] = $null
$data[61,14] = $null
$data[61,15] = $null
$data[61,16] = $null
$data[62,0] = 'East Village'
$data[62,1] = 1
$data[62,2] = $null
$data[62,3] = $null
$data[62,4] = 1
$data[62,5] = $null
$data[62,6] = 1
$data[62,7] = $null
$data[62,8] = $null
$data[62,9] = $null
$data[62,10] = 1
$data[62,11] = $null
$data[62,12] = 1
$data[62,13] = $null
$data[62,14] = $null
$data[62,15] = $null
$data[62,16] = $null
$data[63,0] = 'Fuller Park'
$data[63,1] = $null
$data[63,2] = $null
$data[63,3] = $null
$data[63,4] = $null
$data[63,5] = 1
$data[63,6] = 1
$data[63,7] = $null
$data[63,8] = 1
$data[63,9] = $null
$data[63,10] = $null
$data[63,11] = $null
$data[63,12] = $null
$data[63,13] = $null
$data[63,14] = $null
$data[63,15] = $null
$data[63,16] = $null
$data[64,0] = 'Gage Park'
$data[64,1] = $null
$data[64,2] = $null
$data[64,3] = $null
$data[64,4] = 1
$data[64,5] = 1
$data[64,6] = $null
$data[64,7] = $null
$data[64,8] = 1
$data[64,9] = $null
$data[64,10] = $null
$data[64,11] = $null
$data[64,12] = 2
$data[64,13] = $null
$data[64,14] = 2
$data[64,15] = $null
$data[64,16] = $null
$data[65,0] = 'Garfield Ridge'
$data[65,1] = $null
$data[65,2] = $null
$data[65,3] = $null
$data[65,4] = 1
$data[65,5] = 1
$data[65,6] = $null
$data[65,7] = $null
$data[65,8] = $null
$data[65,9] = $null
$data[65,10] = 1
$data[65,11] = $null
$data[65,12] = $null
$data[65,13] = $null
$data[65,14] = $null
$data[65,15] = $null
$data[65,16] = $null
$data[66,0] = 'Gold Coast'
$data[66,1] = $null
$data[66,2] = $null
$data[66,3] = $null
$data[66,4] = 1
$data[66,5] = $null
$data[66,6] = $null
$data[66,7] = $null
$data[66,8] = $null
$data[66,9] = $null
$data[66,10] = 1
$data[66,11] = $null
$data[66,12] = $null
$data[66,13] = $null
$data[66,14] = $null
$data[66,15] = $null
$data[66,16] = $null
$data[67,0] = 'Hegewisch'
$data[67,1] = $null
$data[67,2] = $null
$data[67,3] = $null
$data[67,4] = $null
$data[67,5] = $null
$data[67,6] = $null
$data[67,7] = $null
$data[67,8] = $null
$data[67,9] = $null
$data[67,10] = $null
$data[67,11] = $null
$data[67,12] = 1
$data[67,13] = $null
$data[67,14] = $null
$data[67,15] = $null
$data[67,16] = $null
$data[68,0] = 'Jackson Park'
$data[68,1] = $null
$data[68,2] = $null
$data[68,3] = $null
$data[68,4] = $null
$data[68,5] = $null
$data[68,6] = $null
$data[68,7] = $null
$data[68,8] = $null
$data[68,9] = $null
$data[68,10] = $null
$data[68,11] = 1
$data[68,12] = 1
$data[68,13] = $null
$data[68,14] = $null
$data[68,15] = $null
$data[68,16] = $null
$data[69,0] = 'Jefferson Park'
$data[69,1] = $null
$data[69,2] = $null
$data[69,3] = $null
$data[69,4] = 3
$data[69,5] = $null
$data[69,6] = $null
$data[69,7] = $null
$data[69,8] = $null
$data[69,9] = $null
$data[69,10] = $null
$data[69,11] = $null
$data[69,12] = $null
$data[69,13] = $null
$data[69,14] = $null
$data[69,15] = $null
$data[69,16] = $null
$data[70,0] = 'Lincoln Park'
$data[70,1] = $null
$data[70,2] = $null
$data[70,3] = 2
$data[70,4] = 2
$data[70,5] = $null
$data[70,6] = $null
$data[70,7] = $null
$data[70,8] = $null
$data[70,9] = $null
$data[70,10] = 1
$data[70,11] = $null
$data[70,12] = $null
$data[70,13] = $null
$data[70,14] = $null
$data[70,15] = $null
$data[70,16] = $null
$data[71,0] = 'Lincoln Square'
$data[71,1] = $null
$data[71,2] = $null
$data[71,3] = $null
$data[71,4] = 1
$data[71,5] = $null
$data[71,6] = $null
$data[71,7] = $null
$data[71,8] = $null
$data[71,9] = $null
$data[71,10] = $null
$data[71,11] = $null
$data[71,12] = $null
$data[71,13] = $null
$data[71,14] = $null
$data[71,15] = $null
$data[71,16] = $null
$data[72,0] = 'Museum Campus'
$data[72,1] = $null
$data[72,2] = $null
$data[72,3] = $null
$data[72,4] = 1
$data[72,5] = $null
$data[72,6] = $null
$data[72,7] = $null
$data[72,8] = $null
$data[72,9] = $null
$data[72,10] = $null
$data[72,11] = $null
$data[72,12] = $null
$data[72,13] = $null
$data[72,14] = $null
$data[72,15] = $null
$data[72,16] = $null
$data[73,0] = 'North Center'
$data[73,1] = $null
$data[73,2] = $null
$data[73,3] = $null
$data[73,4] = $null
$data[73,5] = $null
$data[73,6] = $null
$data[73,7] = $null
$data[73,8] = $null
$data[73,9] = $null
$data[73,10] = $null
$data[73,11] = $null
$data[73,12] = $null
$data[73,13] = $null
$data[73,14] = 1
$data[73,15] = $null
$data[73,16] = $null
$data[74,0] = 'Oakland'
$data[74,1] = $null
$data[74,2] = $null
$data[74,3] = $null
$data[74,4] = 1
$data[74,5] = $null
$data[74,6] = $null
$data[74,7] = $null
$data[74,8] = $null
$data[74,9] = $null
$data[74,10] = $null
$data[74,11] = $null
$data[74,12] = $null
$data[74,13] = $null
$data[74,14] = $null
$data[74,15] = $null
$data[74,16] = $null
$data[75,0] = 'Printers Row'
$data[75,1] = $null
$data[75,2] = $null
$data[75,3] = $null
$data[75,4] = $null
$data[75,5] = $null
$data[75,6] = 1
$data[75,7] = $null
$data[75,8] = $null
$data[75,9] = $null
$data[75,10] = $null
$data[75,11] = $null
$data[75,12] = 1
$data[75,13] = $null
$data[75,14] = $null
$data[75,15] = $null
$data[75,16] = $null
$data[76,0] = 'Pullman'
$data[76,1] = $null
$data[76,2] = $null
$data[76,3] = $null
$data[76,4] = $null
$data[76,5] = $null
$data[76,6] = $null
$data[76,7] = $null
$data[76,8] = $null
$data[76,9] = $null
$data[76,10] = $null
$data[76,11] = $null
$data[76,12] = $null
$data[76,13] = $null
$data[76,14] = $null
$data[76,15] = 1
$data[76,16] = $null
$data[77,0] = 'Rush & Division'
$data[77,1] = $null
$data[77,2] = $null
$data[77,3] = $null
$data[77,4] = $null
$data[77,5] = 1
$data[77,6] = $null
$data[77,7] = $null
$data[77,8] = $null
$data[77,9] = 1
$data[77,10] = $null
$data[77,11] = $null
$data[77,12] = $null
$data[77,13] = $null
$data[77,14] = $null
$data[77,15] = $null
$data[77,16] = $null
$data[78,0] = 'Sauganash,Forest Glen'
$data[78,1] = $null
$data[78,2] = $null
$data[78,3] = $null
$data[78,4] = 2
$data[78,5] = $null
$data[78,6] = $null
$data[78,7] = $null
$data[78,8] = $null
$data[78,9] = $null
$data[78,10] = $null
$data[78,11] = $null
$data[78,12] = $null
$data[78,13] = $null
$data[78,14] = $null
$data[78,15] = $null
$data[78,16] = $null
$data[79,0] = 'Sheffield & DePaul'
$data[79,1] = $null
$data[79,2] = $null
$data[79,3] = $null
$data[79,4] = 1
$data[79,5] = $null
$data[79,6] = $null
$data[79,7] = $null
$data[79,8] = $null
$data[79,9] = $null
$data[79,10] = $null
$data[79,11] = $null
$data[79,12] = $null
$data[79,13] = $null
$data[79,14] = $null
$data[79,15] = $null
$data[79,16] = $null
$data[80,0] = 'Streeterville'
$data[80,1] = $null
$data[80,2] = $null
$data[80,3] = $null
$data[80,4] = $null
$data[80,5] = $null
$data[80,6] = $null
$data[80,7] = $null
$data[80,8] = $null
$data[80,9] = 1
$data[80,10] = $null
$data[80,11] = $null
$data[80,12] = $null
$data[80,13] = $null
$data[80,14] = $null
$data[80,15] = $null
$data[80,16] = $null
$data[81,0] = 'Ukrainian Village'
$data[81,1] = $null
$data[81,2] = $null
$data[81,3] = 3
$data[81,4] = 3
$data[81,5] = 1
$data[81,6] = 1
$data[81,7] = $null
$data[81,8] = $null
$data[81,9] = 1
$data[81,10] = $null
$data[81,11] = 1
$data[81,12] = $null
$data[81,13] = $null
$data[81,14] = $null
$data[81,15] = $null
$data[81,16] = $null
$data[82,0] = 'Washington Park'
$data[82,1] = 1
$data[82,2] = $null
$data[82,3] = $null
$data[82,4] = $null
$data[82,5] = $null
$data[82,6] = 2
$data[82,7] = 1
$data[82,8] = $null
$data[82,9] = $null
$data[82,10] = 4
$data[82,11] = 1
$data[82,12] = $null
$data[82,13] = $null
$data[82,14] = 1
$data[82,15] = $null
$data[82,16] = $null
$data[83,0] = 'West Elsdon'
$data[83,1] = $null
$data[83,2] = $null
$data[83,3] = $null
$data[83,4] = $null
$data[83,5] = $null
$data[83,6] = $null
$data[83,7] = $null
$data[83,8] = $null
$data[83,9] = 1
$data[83,10] = $null
$data[83,11] = $null
$data[83,12] = 1
$data[83,13] = $null
$data[83,14] = $null
$data[83,15] = $null
$data[83,16] = $null
$data[84,0] = 'Wrigleyville'
$data[84,1] = $null
$data[84,2] = $null
$data[84,3] = $null
$data[84,4] = $null
$data[84,5] = $null
$data[84,6] = $null
$data[84,7] = $null
$data[84,8] = $null
$data[84,9] = $null
$data[84,10] = $null
$data[84,11] = $null
$data[84,12] = 1
$data[84,13] = $null
$data[84,14] = $null
$data[84,15] = $null
$data[84,16] = $null

$ws.Range("A2:Q86").Value = $data

